# VMS_front_rear.xlsx edit:
#  - Update point "G" (row 8) X/Y/Z coordinates on both the front and the
#    rear suspension sheets.
#  - Move the active-cell selection from B17 to B9 on both sheets (front
#    sheet stays the active tab).
#  - Slightly resize/reposition the reference pictures on both sheets to
#    match the new drawing geometry.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # VMS Front Suspension
$ws2 = $wb.Worksheets.Item(2)   # VMS Rear Suspension

# --- Point G coordinates (row 8) ------------------------------------------
$ws1.Range("B8").Value = 585.95
$ws1.Range("C8").Value = 287.25
$ws1.Range("D8").Value = 1520.13

$ws2.Range("B8").Value = 616.11
$ws2.Range("C8").Value = 153.79
$ws2.Range("D8").Value = -8.21

# --- Picture geometry (keep the same top-left anchor point, tweak size) ---
$shp1 = $ws1.Shapes.Item(1)
$shp1.Left   = 184.79055118110236
$shp1.Top    = 11.537007874015748
$shp1.Width  = 299.36692913385826
$shp1.Height = 367.5968503937008

$shp2 = $ws2.Shapes.Item(1)
$shp2.Left   = 206.2204724409449
$shp2.Top    = 15.73228346456693
$shp2.Width  = 343.7574803149606
$shp2.Height = 422.0220472440945

# --- Selection: move from B17 to B9 on both sheets, front sheet active ----
[void]$ws2.Range("B9").Select()
[void]$ws1.Range("B9").Select()
